# Text updates as supplied by PM&C.
# Applies to the "Description" worksheet of the Health - Type 2 Diabetes
# dashboard workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# --- B7 ("Influences" body text) -------------------------------------------------
# Drop the quotation marks / citation and reword "overweight or obesity" ->
# "being overweight or obese".
$b7Font = $ws.Range("B7").Font
$b7Font.Name = "Calibri"
$b7Font.Size = 11
$ws.Range("B7").WrapText = $true
$ws.Range("B7").Value = "The prevalence of Type 2 diabetes is increased through hereditary factors and lifestyle risk factors including poor diet, insufficient physical activity and being overweight or obese."
$ws.Rows.Item(7).RowHeight = 25.45

# --- A11 / B11 (new "Source" label + shortened source text) ----------------------
# Previously A11 was blank and B11 held "Sourced from: ...". Now the row gets an
# explicit "Source" label in column A (matching the other A-column labels) and the
# "Sourced from: " prefix is dropped from the citation text itself.
$ws.Range("A11").Value = "Source"

$b11Font = $ws.Range("B11").Font
$b11Font.Name = "Calibri"
$b11Font.Size = 11
$ws.Range("B11").WrapText = $true
$ws.Range("B11").Value = "ABS unpublished, Australian Health Survey 2011–13 (2011-12 NHMS component); ABS unpublished, Australian Aboriginal and Torres Strait Islander Health Survey, 2012-13 (National Aboriginal and Torres Strait Islander Health Measures Survey component)."
$ws.Rows.Item(11).RowHeight = 37.45

# Keep the active selection on B7, matching the re-saved workbook.
$ws.Range("B7").Select()
